$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (A1:D1) and retitle municipality/state text cells ---
# (Spanish connector words de/del/el/la/los/las/y -> title case De/Del/El/La/Los/Las/Y,
#  plus a couple of literal corrections e.g. MonteMorelos -> Montemorelos)
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'
$ws.Range('B5').Value = 'Pabellón De Arteaga'
$ws.Range('B6').Value = 'Rincón De Romos'
$ws.Range('B19').Value = 'Amatenango De La Frontera'
$ws.Range('B25').Value = 'Chiapa De Corzo'
$ws.Range('B28').Value = 'Comitán De Domínguez'
$ws.Range('B44').Value = 'Marqués De Comillas'
$ws.Range('B45').Value = 'Mazapa De Madero'
$ws.Range('B50').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B55').Value = 'San Cristóbal De Las Casas'
$ws.Range('B80').Value = 'Guadalupe Y Calvo'
$ws.Range('B82').Value = 'Hidalgo Del Parral'
$ws.Range('B97').Value = 'San Juan De Sabinas'
$ws.Range('A107').Value = 'Ciudad De México'
$ws.Range('B111').Value = 'Cuajimalpa De Morelos'
$ws.Range('B138').Value = 'Pánuco De Coronado'
$ws.Range('A146').Value = 'Estado De México'
$ws.Range('B146').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B148').Value = 'Almoloya De Juárez'
$ws.Range('B150').Value = 'Atizapán De Zaragoza'
$ws.Range('B152').Value = 'Chapa De Mota'
$ws.Range('B156').Value = 'Ecatepec De Morelos'
$ws.Range('B159').Value = 'Ixtapan De La Sal'
$ws.Range('B168').Value = 'Naucalpan De Juárez'
$ws.Range('B173').Value = 'San Felipe Del Progreso'
$ws.Range('B179').Value = 'Tenango Del Valle'
$ws.Range('B185').Value = 'Tlalnepantla De Baz'
$ws.Range('B189').Value = 'Valle De Bravo'
$ws.Range('B198').Value = 'Apaseo El Alto'
$ws.Range('B199').Value = 'Apaseo El Grande'
$ws.Range('B204').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B207').Value = 'Jaral Del Progreso'
$ws.Range('B216').Value = 'San Diego De La Unión'
$ws.Range('B218').Value = 'San Francisco Del Rincón'
$ws.Range('B220').Value = 'San Luis De La Paz'
$ws.Range('B221').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B223').Value = 'Silao De La Victoria'
$ws.Range('B228').Value = 'Valle De Santiago'
$ws.Range('B234').Value = 'Acapulco De Juárez'
$ws.Range('B237').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B241').Value = 'Atenango Del Río'
$ws.Range('B243').Value = 'Atoyac De Álvarez'
$ws.Range('B244').Value = 'Ayutla De Los Libres'
$ws.Range('B247').Value = 'Chilapa De Álvarez'
$ws.Range('B248').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B251').Value = 'Coyuca De Benítez'
$ws.Range('B252').Value = 'Coyuca De Catalán'
$ws.Range('B255').Value = 'Cutzamala De Pinzón'
$ws.Range('B259').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B260').Value = 'Iguala De La Independencia'
$ws.Range('B261').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B262').Value = 'Zihuatanejo De Azueta'
$ws.Range('B264').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B265').Value = 'Mártir De Cuilapan'
$ws.Range('B276').Value = 'Taxco De Alarcón'
$ws.Range('B278').Value = 'Técpan De Galeana'
$ws.Range('B280').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B282').Value = 'Tixtla De Guerrero'
$ws.Range('B293').Value = 'Agua Blanca De Iturbide'
$ws.Range('B297').Value = 'Atotonilco El Grande'
$ws.Range('B300').Value = 'Cuautepec De Hinojosa'
$ws.Range('B303').Value = 'Huasca De Ocampo'
$ws.Range('B309').Value = 'Molango De Escamilla'
$ws.Range('B311').Value = 'Nopala De Villagrán'
$ws.Range('B312').Value = 'Pachuca De Soto'
$ws.Range('B316').Value = 'Santiago De Anaya'
$ws.Range('B320').Value = 'Tenango De Doria'
$ws.Range('B322').Value = 'Tepehuacán De Guerrero'
$ws.Range('B323').Value = 'Tezontepec De Aldama'
$ws.Range('B328').Value = 'Tulancingo De Bravo'
$ws.Range('B330').Value = 'Zacualtipán De Ángeles'
$ws.Range('B333').Value = 'Ahualulco De Mercado'
$ws.Range('B336').Value = 'Atotonilco El Alto'
$ws.Range('B346').Value = 'Encarnación De Díaz'
$ws.Range('B349').Value = 'Ixtlahuacán Del Río'
$ws.Range('B353').Value = 'Lagos De Moreno'
$ws.Range('B357').Value = 'Ojuelos De Jalisco'
$ws.Range('B360').Value = 'San Juan De Los Lagos'
$ws.Range('B363').Value = 'San Miguel El Alto'
$ws.Range('B364').Value = 'San Sebastián Del Oeste'
$ws.Range('B366').Value = 'Tamazula De Gordiano'
$ws.Range('B370').Value = 'Tepatitlán De Morelos'
$ws.Range('B371').Value = 'Tizapán El Alto'
$ws.Range('B372').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B378').Value = 'Unión De San Antonio'
$ws.Range('B379').Value = 'Valle De Juárez'
$ws.Range('B381').Value = 'Yahualica De González Gallo'
$ws.Range('B382').Value = 'Zacoalco De Torres'
$ws.Range('B384').Value = 'Zapotlán El Grande'
$ws.Range('B398').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B443').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B463').Value = 'Coatlán Del Río'
$ws.Range('B471').Value = 'Puente De Ixtla'
$ws.Range('B475').Value = 'Tetela Del Volcán'
$ws.Range('B476').Value = 'Tlaltizapán De Zapata'
$ws.Range('B481').Value = 'Zacualpan De Amilpas'
$ws.Range('B501').Value = 'Mier Y Noriega'
$ws.Range('B502').Value = 'Montemorelos'
$ws.Range('B505').Value = 'San Nicolás De Los Garza'
$ws.Range('B510').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B513').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B516').Value = 'Guelatao De Juárez'
$ws.Range('B517').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B518').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B519').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B520').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B525').Value = 'Mariscala De Juárez'
$ws.Range('B528').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B530').Value = 'Oaxaca De Juárez'
$ws.Range('B531').Value = 'Ocotlán De Morelos'
$ws.Range('B532').Value = 'Putla Villa De Guerrero'
$ws.Range('B535').Value = 'San Antonino El Alto'
$ws.Range('B555').Value = 'San Miguel El Grande'
$ws.Range('B595').Value = 'Santo Domingo De Morelos'
$ws.Range('B601').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B602').Value = 'Teotitlán De Flores Magón'
$ws.Range('B604').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B605').Value = 'Tlacolula De Matamoros'
$ws.Range('B606').Value = 'Tlalixtac De Cabrera'
$ws.Range('B607').Value = 'Totontepec Villa De Morelos'
$ws.Range('B608').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B609').Value = 'Villa De Zaachila'
$ws.Range('B610').Value = 'Zimatlán De Álvarez'
$ws.Range('B619').Value = 'Chalchicomula De Sesma'
$ws.Range('B632').Value = 'Huehuetlán El Chico'
$ws.Range('B635').Value = 'Izúcar De Matamoros'
$ws.Range('B638').Value = 'Los Reyes De Juárez'
$ws.Range('B657').Value = 'Tepexi De Rodríguez'
$ws.Range('B658').Value = 'Tetela De Ocampo'
$ws.Range('B661').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B667').Value = 'Xayacatlán De Bravo'
$ws.Range('B677').Value = 'Jalpan De Serra'
$ws.Range('B679').Value = 'Pinal De Amoles'
$ws.Range('B691').Value = 'Axtla De Terrazas'
$ws.Range('B702').Value = 'San Ciro De Acosta'
$ws.Range('B706').Value = 'Santa María Del Río'
$ws.Range('B707').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B711').Value = 'Tanquián De Escobedo'
$ws.Range('B714').Value = 'Villa De Arista'
$ws.Range('B715').Value = 'Villa De La Paz'
$ws.Range('B716').Value = 'Villa De Ramos'
$ws.Range('B717').Value = 'Villa De Reyes'
$ws.Range('B777').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B788').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B795').Value = 'Boca Del Río'
$ws.Range('B796').Value = 'Camarón De Tejeda'
$ws.Range('B798').Value = 'Castillo De Teayo'
$ws.Range('B804').Value = 'Cosamaloapan De Carpio'
$ws.Range('B811').Value = 'Hueyapan De Ocampo'
$ws.Range('B812').Value = 'Ignacio De La Llave'
$ws.Range('B816').Value = 'Ixhuacán De Los Reyes'
$ws.Range('B822').Value = 'Juchique De Ferrer'
$ws.Range('B825').Value = 'Lerdo De Tejada'
$ws.Range('B826').Value = 'Martínez De La Torre'
$ws.Range('B836').Value = 'Paso De Ovejas'
$ws.Range('B837').Value = 'Paso Del Macho'
$ws.Range('B839').Value = 'Poza Rica De Hidalgo'
$ws.Range('B842').Value = 'Soledad De Doblado'
$ws.Range('B846').Value = 'Tatahuicapan De Juárez'
$ws.Range('B862').Value = 'Vega De Alatorre'
$ws.Range('B873').Value = 'Concepción Del Oro'
$ws.Range('B883').Value = 'Jiménez Del Teul'
$ws.Range('B890').Value = 'Moyahua De Estrada'
$ws.Range('B891').Value = 'Nochistlán De Mejía'
$ws.Range('B892').Value = 'Noria De Ángeles'
$ws.Range('B901').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B903').Value = 'Villa De Cos'

# --- Refresh percentage column precision (re-store the computed share so it round-trips
#     through Excel's 15-significant-digit storage, matching the canonical file) ---
$ws.Range('D9').Value2 = 0.000952154248988336
$ws.Range('D26').Value2 = 0.000952154248988336
$ws.Range('D30').Value2 = 0.000952154248988336
$ws.Range('D39').Value2 = 0.000952154248988336
$ws.Range('D66').Value2 = 0.000952154248988336
$ws.Range('D68').Value2 = 0.000952154248988336
$ws.Range('D71').Value2 = 0.000952154248988336
$ws.Range('D78').Value2 = 0.000952154248988336
$ws.Range('D118').Value2 = 0.000952154248988336
$ws.Range('D137').Value2 = 0.000952154248988336
$ws.Range('D187').Value2 = 0.000952154248988336
$ws.Range('D223').Value2 = 0.000952154248988336
$ws.Range('D253').Value2 = 0.000952154248988336
$ws.Range('D258').Value2 = 0.000952154248988336
$ws.Range('D260').Value2 = 0.000952154248988336
$ws.Range('D263').Value2 = 0.000952154248988336
$ws.Range('D266').Value2 = 0.000952154248988336
$ws.Range('D270').Value2 = 0.000952154248988336
$ws.Range('D277').Value2 = 0.000952154248988336
$ws.Range('D292').Value2 = 0.000952154248988336
$ws.Range('D298').Value2 = 0.000952154248988336
$ws.Range('D314').Value2 = 0.000952154248988336
$ws.Range('D322').Value2 = 0.000952154248988336
$ws.Range('D369').Value2 = 0.000952154248988336
$ws.Range('D384').Value2 = 0.000952154248988336
$ws.Range('D407').Value2 = 0.000952154248988336
$ws.Range('D432').Value2 = 0.000952154248988336
$ws.Range('D441').Value2 = 0.000952154248988336
$ws.Range('D444').Value2 = 0.000952154248988336
$ws.Range('D451').Value2 = 0.000952154248988336
$ws.Range('D471').Value2 = 0.000952154248988336
$ws.Range('D520').Value2 = 0.000952154248988336
$ws.Range('D525').Value2 = 0.000952154248988336
$ws.Range('D546').Value2 = 0.000952154248988336
$ws.Range('D577').Value2 = 0.000952154248988336
$ws.Range('D581').Value2 = 0.000952154248988336
$ws.Range('D585').Value2 = 0.000952154248988336
$ws.Range('D594').Value2 = 0.000952154248988336
$ws.Range('D608').Value2 = 0.000952154248988336
$ws.Range('D610').Value2 = 0.000952154248988336
$ws.Range('D613').Value2 = 0.000952154248988336
$ws.Range('D657').Value2 = 0.000952154248988336
$ws.Range('D667').Value2 = 0.000952154248988336
$ws.Range('D668').Value2 = 0.000952154248988336
$ws.Range('D687').Value2 = 0.000952154248988336
$ws.Range('D708').Value2 = 0.000952154248988336
$ws.Range('D735').Value2 = 0.000952154248988336
$ws.Range('D736').Value2 = 0.000952154248988336
$ws.Range('D739').Value2 = 0.000952154248988336
$ws.Range('D761').Value2 = 0.000952154248988336
$ws.Range('D770').Value2 = 0.000952154248988336
$ws.Range('D771').Value2 = 0.000952154248988336
$ws.Range('D807').Value2 = 0.000952154248988336
$ws.Range('D809').Value2 = 0.000952154248988336
$ws.Range('D831').Value2 = 0.000952154248988336
$ws.Range('D858').Value2 = 0.000952154248988336
$ws.Range('D860').Value2 = 0.000952154248988336
$ws.Range('D880').Value2 = 0.000952154248988336
$ws.Range('D886').Value2 = 0.000952154248988336
$ws.Range('D901').Value2 = 0.000952154248988336
$ws.Range('D903').Value2 = 0.000952154248988336
$ws.Range('D905').Value2 = 0.000952154248988336

# --- Drop the trailing footnote/metadata block (rows 911-915); row 910 was already blank ---
$ws.Rows("911:915").Delete()

